$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.444093465805054
$ws.Range("B1").Value = 1.972862362861633
$ws.Range("C1").Value = 3.044617414474487
$ws.Range("D1").Value = 4.888361930847168
$ws.Range("E1").Value = 0.9191415905952454
